# Update Name of Algo
# Apply updated values to the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = -10.92689999999999

# Row 4
$ws.Range("A4").Value = -21.14310000000001
$ws.Range("B4").Value = 4.837100000000005
$ws.Range("C4").Value = -11.04319999999999

# Row 5
$ws.Range("B5").Value = 5.131400000000001

# Row 6
$ws.Range("A6").Value = -21.3928

# Row 7
$ws.Range("A7").Value = -21.4279

# Row 8
$ws.Range("B8").Value = 4.861400000000003

# Row 9
$ws.Range("C9").Value = -11.81610000000001

# Row 11
$ws.Range("C11").Value = -14.0186

# Row 14
$ws.Range("C14").Value = -11.9723

# Row 16
$ws.Range("A16").Value = -21.46170000000002
$ws.Range("B16").Value = 5.439699999999995

# Row 18
$ws.Range("C18").Value = -14.50530000000001

# Row 20
$ws.Range("A20").Value = -22.89900000000002

# Row 22
$ws.Range("B22").Value = 5.4118

# Row 25
$ws.Range("C25").Value = -10.91529999999999
